$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 312 (shifts existing rows 312..376 down to 313..377)
$ws.Rows.Item(312).Insert()

# Populate the newly inserted row 312 with the new data record
$ws.Cells.Item(312, 1).Value = 3
$ws.Cells.Item(312, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(312, 3).Value = "Coquimbo"
$ws.Cells.Item(312, 4).Value = 44711
$ws.Cells.Item(312, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(312, 5).Value = 5
$ws.Cells.Item(312, 6).Value = 100112031
$ws.Cells.Item(312, 7).Value = "Poroto verde"
$ws.Cells.Item(312, 8).Value = "Magnum"
$ws.Cells.Item(312, 9).Value = "Primera"
$ws.Cells.Item(312, 10).Value = 76
$ws.Cells.Item(312, 11).Value = 27000
$ws.Cells.Item(312, 12).Value = 28000
$ws.Cells.Item(312, 13).Value = 27500
$ws.Cells.Item(312, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(312, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(312, 16).Value = 1100
$ws.Cells.Item(312, 17).Value = 25
$ws.Cells.Item(312, 18).Value = "Hortaliza"
